$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

$c2text = @'
### Dirty install
* Make sure that you have several versions of Mantid installed
 * Last release
 * A nightly
 * If possible an old release
* Install the new release candidate of Mantid
- [ ] Check that Mantid boots up correctly
### Clean install
* Remove all existing Mantid versions and associated files
**On Windows**:
 * Uninstall the program
 * Clear shortcuts from desktop
 * Remove the mantid  files in %APPDATA%
**On macOS** :
 * Remove the application
 * Remove the `~/.mantid directory`
 * Remove (or rename) `~/Library/Preferences/org.python.python.Python.plist`  and `~/Library/Preferences/org.python.python.plist` 
**On Linux** :
 * Remove the package: `/opt/Mantid`
 * Remove `~/.config/Mantid`
 * Remove `~/.mantid/`
Re-install the new release candidate of Mantid
- [ ] Check that Mantid boots up correctly

'@
$c4text = @'
## ISIS only, if possible, so catalogue access can be tested
- [ ] MantidWorkbench opens without errors or warnings 
- [ ] Every option in `Interface` opens a GUI 
- [ ] Load some test data 
- [ ] Access Catalogue through algorithms (use Facilities account, same as for IDAaaS): [CatalogLogin](https://docs.mantidproject.org/algorithms/CatalogLogin-v1.html) and [CatalogGetDataFiles](https://docs.mantidproject.org/algorithms/CatalogGetDataFiles-v1.html)
- [ ] Saving/loading projects works 
 - [ ] Alter preferences in [File > Settings](https://docs.mantidproject.org/nightly/tutorials/mantid_basic_course/useful/01_save_settings.html#settings) and check they are obeyed
'@
$c5text = @'
* 1D plotting:[instructions](http://docs.mantidproject.org/nightly/tutorials/mantid_basic_course/loading_and_displaying_data/03_displaying_1D_data.html#displaying-1d-data)
 - [ ] Simple plot 
 - [ ] Another way to plot 
 - [ ] Adding curves to existing plots 
 - [ ] Also, test out [waterfall](https://docs.mantidproject.org/nightly/plotting/WaterfallPlotsHelp.html#waterfall-plots) and [tiled]( https://docs.mantidproject.org/nightly/plotting/1DPlotsHelp.html#tiled-plots)
 - [ ] Check Toolbar buttons

* 2D plotting: [instructions](http://docs.mantidproject.org/nightly/tutorials/mantid_basic_course/loading_and_displaying_data/04_displaying_2D_data.html#displaying-2d-data)
 - [ ] Plot all spectra 
 - [ ] Change colour map 
 - [ ] [Contour plot](https://docs.mantidproject.org/nightly/plotting/ColorfillPlotsHelp.html#contour-plots) (under 3D menu) 
 - [ ] Check Toolbar buttons

* [3D plotting](https://docs.mantidproject.org/nightly/plotting/3DPlotsHelp.html):
 - Load some data eg `LOQ7041` from the ISIS sample data
 - [ ] 3D surface
 - [ ] 3D wire frame  
 - [ ] Check Toolbar buttons

 ## Sliceviewer
 - [ ] Overly long instructions (don't spend  long!) and data [here](https://developer.mantidproject.org/Testing/SliceViewer/SliceViewer.html). In particular try editing the data in a workspace while it is open in Sliceviewer!
'@
$c6text = @'
* Test that the Python scripting window works, [directions here](https://docs.mantidproject.org/nightly/workbench/scriptwindow.html)
 - [ ] Editor options 
 - [ ] Execution options 
 - [ ] Script output  
- [ ] Perform some workspace algebra 
- [ ] Test numpy functionality 
- [ ] Use the scripting window to run some scripts 
- [ ] Run through some examples from [the documentation](https://docs.mantidproject.org/nightly/tutorials/python_in_mantid/index.html), 3 or 4 examples from the Solutions is enough 
'@
$condaDescText = @'
Quickly run through some of the other Smoke testing instructions on the separate Conda release package. The most useful tests is to check many different dependencies, such as numpy and matpltlib in the editor, and opening the interfaces.
'@

$ws.Range("C2").Value2 = $c2text
$ws.Range("C4").Value2 = $c4text
$ws.Range("C5").Value2 = $c5text
$ws.Range("C6").Value2 = $c6text

$ws.Range("A8").Value2 = "Conda Package Tests"
$ws.Range("B8").Value2 = " :snake:"
$ws.Range("C8").Value2 = $condaDescText
$ws.Rows.Item(8).RowHeight = 48

# Restore the original row heights for rows 2-7: updating the wrapped text in
# column C can make this headless host's approximate autofit logic resize
# these rows, but the real workbook keeps their original custom heights.
$ws.Rows.Item(2).RowHeight = 135
$ws.Rows.Item(3).RowHeight = 90.75
$ws.Rows.Item(4).RowHeight = 74.25
$ws.Rows.Item(5).RowHeight = 92.25
$ws.Rows.Item(6).RowHeight = 96.75
$ws.Rows.Item(7).RowHeight = 99.75

$ws.Range("C9").Select() | Out-Null

